# "Generate Report for Handback" — refresh the handoff/handback timestamp
# columns that get re-stamped whenever the handback status report is
# regenerated.
#
# Overview sheet, column G = "Latest HO Xliff Generate Date"
# zh-cn / de-de sheets, column H = "Correspond Handoff Datetime"
# zh-cn / de-de sheets, column K = "Correspond Handback DateTime"
#
# Only row 2 (the 68b2042d-... file) changes in this run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$dateFormat = "yyyy-mm-dd HH:mm:ss"

# Overview!G2 — Latest HO Xliff Generate Date
# (de-de!H2 "Correspond Handoff Datetime" carried this exact same
#  timestamp text and must be refreshed together with it.)
$wsOverview.Range("G2").Value = "2016-09-01 17:15:21"
$wsOverview.Range("G2").NumberFormat = $dateFormat

$wsDeDe.Range("H2").Value = "2016-09-01 17:15:21"
$wsDeDe.Range("H2").NumberFormat = $dateFormat

# zh-cn!H2 — Correspond Handoff Datetime
$wsZhCn.Range("H2").Value = "2016-09-01 17:15:10"
$wsZhCn.Range("H2").NumberFormat = $dateFormat

# zh-cn!K2 — Correspond Handback DateTime
$wsZhCn.Range("K2").Value = "2016-09-01 17:15:57"
$wsZhCn.Range("K2").NumberFormat = $dateFormat

# de-de!K2 — Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-09-01 17:16:13"
$wsDeDe.Range("K2").NumberFormat = $dateFormat
